# 1. Remove the stray empty <c> placeholders in column B of "ODI Batting"
#    for rows 3, 4, 7, 8 and 10 (they hold no real value; the authoritative
#    edit simply drops the inline-string cell entirely).
$wb = $excel.ActiveWorkbook
$batting = $wb.Worksheets.Item("ODI Batting")
foreach ($r in 3, 4, 7, 8, 10) {
    $batting.Cells.Item($r, 2).ClearContents()
}

# 2. Append a new worksheet "ODI Batting Extra" as the 4th/last sheet and
#    populate it with the MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#    PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH table.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Header row - copy the header formatting (bold, centered, thin border)
# from the existing "ODI Batting" header so the new sheet's header reuses
# the same visual style as the rest of the workbook.
$extra.Cells.Item(1, 1).Value = "MATCH_CODE"
$extra.Cells.Item(1, 2).Value = "BATTING_POSITION"
$extra.Cells.Item(1, 3).Value = "NUM_4"
$extra.Cells.Item(1, 4).Value = "NUM_6"
$extra.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

$batting.Range("A1:F1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)

# Data rows. MATCH_CODE values look numeric but are authored as text, and
# the percentage column is literal text too ("2.23%"), so those go in with
# a leading apostrophe to force text storage instead of Excel's automatic
# number / percent coercion. BATTING_POSITION is a genuine number.
$data = @(
    @("4564", 9,  "0", "0", "",       "NO"),
    @("4565", $null, "", "", "",     "NO"),
    @("4567", 8,  "", "", "",         "NO"),
    @("4590", $null, "", "", "",     "NO"),
    @("4592", 9,  "1", "0", "2.23%",  "NO"),
    @("4634", $null, "", "", "",     "NO"),
    @("4638", 9,  "", "", "",         "NO"),
    @("4641", 8,  "0", "1", "5.34%",  "NO"),
    @("4686", $null, "", "", "",     "NO"),
    @("4688", 9,  "1", "0", "5.49%",  "NO"),
    @("4690", $null, "", "", "",     "NO")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]

    $extra.Cells.Item($row, 1).Value = "'" + $vals[0]

    if ($null -ne $vals[1]) {
        $extra.Cells.Item($row, 2).Value = $vals[1]
    }

    if ($vals[2] -ne "") {
        $extra.Cells.Item($row, 3).Value = "'" + $vals[2]
    }

    if ($vals[3] -ne "") {
        $extra.Cells.Item($row, 4).Value = "'" + $vals[3]
    }

    if ($vals[4] -ne "") {
        $extra.Cells.Item($row, 5).Value = "'" + $vals[4]
    }

    $extra.Cells.Item($row, 6).Value = $vals[5]
}
